$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $style = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $style
}

Set-TextValue $ws.Range("D2") "30.196.35"
Set-TextValue $ws.Range("E2") "  -3.36%  "
Set-TextValue $ws.Range("D3") "1.922.22"
Set-TextValue $ws.Range("E3") "  -3.02%  "
Set-TextValue $ws.Range("D4") "0.9994"
Set-TextValue $ws.Range("E4") "  -0.07%  "
Set-TextValue $ws.Range("D5") "246.42"
Set-TextValue $ws.Range("E5") "  -2.54%  "
Set-TextValue $ws.Range("D6") "0.7184"
Set-TextValue $ws.Range("E6") "  -9.61%  "
Set-TextValue $ws.Range("D7") "0.9994"
Set-TextValue $ws.Range("E7") "  -0.09%  "
Set-TextValue $ws.Range("D8") "0.3244"
Set-TextValue $ws.Range("E8") "  -3.85%  "
Set-TextValue $ws.Range("D9") "26.28"
Set-TextValue $ws.Range("E9") "  +2.57%  "
Set-TextValue $ws.Range("D10") "0.06841"
Set-TextValue $ws.Range("E10") "  -1.11%  "
Set-TextValue $ws.Range("D11") "0.7940"
Set-TextValue $ws.Range("E11") "  -4.66%  "
Set-TextValue $ws.Range("D12") "0.07921"
Set-TextValue $ws.Range("E12") "  -2.17%  "
Set-TextValue $ws.Range("D13") "1.919.25"
Set-TextValue $ws.Range("E13") "  -3.21%  "
Set-TextValue $ws.Range("D14") "5.385"
Set-TextValue $ws.Range("E14") "  -1.09%  "
Set-TextValue $ws.Range("D15") "94.45"
Set-TextValue $ws.Range("E15") "  -5.58%  "
Set-TextValue $ws.Range("D16") "14.43"
Set-TextValue $ws.Range("E16") "  +4.48%  "
Set-TextValue $ws.Range("D17") "258.90"
Set-TextValue $ws.Range("E17") "  -5.20%  "
Set-TextValue $ws.Range("D18") "30.201.23"
Set-TextValue $ws.Range("E18") "  -3.32%  "
Set-TextValue $ws.Range("D19") "5.827"
Set-TextValue $ws.Range("E19") "  +2.12%  "
Set-TextValue $ws.Range("D20") "0.000007907"
Set-TextValue $ws.Range("E20") "  -0.02%  "
Set-TextValue $ws.Range("D21") "2.169.36"
Set-TextValue $ws.Range("E21") "  -3.36%  "
Set-TextValue $ws.Range("D22") "0.9989"
Set-TextValue $ws.Range("E22") "  -0.34%  "
Set-TextValue $ws.Range("D23") "0.9983"
Set-TextValue $ws.Range("E23") "  -0.28%  "
Set-TextValue $ws.Range("D24") "6.848"
Set-TextValue $ws.Range("E24") "  -1.09%  "
Set-TextValue $ws.Range("D25") "9.675"
Set-TextValue $ws.Range("E25") "  +0.65%  "
Set-TextValue $ws.Range("D26") "160.53"
Set-TextValue $ws.Range("E26") "  -2.39%  "
Set-TextValue $ws.Range("D27") "0.1334"
Set-TextValue $ws.Range("E27") "  -9.62%  "
Set-TextValue $ws.Range("D28") "18.78"
Set-TextValue $ws.Range("E28") "  -4.85%  "
Set-TextValue $ws.Range("D29") "2.236"
Set-TextValue $ws.Range("E29") "  +2.73%  "
Set-TextValue $ws.Range("D30") "1.358"
Set-TextValue $ws.Range("E30") "  +0.14%  "
Set-TextValue $ws.Range("D31") "1.545"
Set-TextValue $ws.Range("E31") "  -1.24%  "
Set-TextValue $ws.Range("D32") "4.408"
Set-TextValue $ws.Range("E32") "  -3.14%  "
Set-TextValue $ws.Range("D33") "4.184"
Set-TextValue $ws.Range("E33") "  -3.15%  "
Set-TextValue $ws.Range("D34") "0.05031"
Set-TextValue $ws.Range("E34") "  -2.22%  "
Set-TextValue $ws.Range("D35") "1.191"
Set-TextValue $ws.Range("E35") "  -1.23%  "
Set-TextValue $ws.Range("D36") "0.7391"
Set-TextValue $ws.Range("E36") "  -1.80%  "
Set-TextValue $ws.Range("D37") "2.730"
Set-TextValue $ws.Range("E37") "  -2.10%  "
Set-TextValue $ws.Range("D38") "0.01944"
Set-TextValue $ws.Range("E38") "  -2.87%  "
Set-TextValue $ws.Range("D39") "2.804"
Set-TextValue $ws.Range("E39") "  -3.62%  "
Set-TextValue $ws.Range("D40") "79.92"
Set-TextValue $ws.Range("E40") "  +2.35%  "
Set-TextValue $ws.Range("D41") "6.482"
Set-TextValue $ws.Range("E41") "  -1.72%  "
Set-TextValue $ws.Range("D42") "0.4417"
Set-TextValue $ws.Range("E42") "  -4.65%  "
Set-TextValue $ws.Range("D43") "2.005"
Set-TextValue $ws.Range("E43") "  -2.16%  "
Set-TextValue $ws.Range("D44") "0.9998"
Set-TextValue $ws.Range("E44") "  -0.11%  "
Set-TextValue $ws.Range("D45") "0.8306"
Set-TextValue $ws.Range("E45") "  -2.34%  "
Set-TextValue $ws.Range("D46") "102.20"
Set-TextValue $ws.Range("E46") "  -2.85%  "
Set-TextValue $ws.Range("D47") "9.683"
Set-TextValue $ws.Range("E47") "  -3.11%  "
Set-TextValue $ws.Range("D48") "7.259"
Set-TextValue $ws.Range("E48") "  -2.74%  "
Set-TextValue $ws.Range("E49") "  -1.12%  "
Set-TextValue $ws.Range("D50") "0.4096"
Set-TextValue $ws.Range("E50") "  -3.85%  "
Set-TextValue $ws.Range("D51") "1.471"
